$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3344
$ws1.Range("G3").Value = 80
$ws1.Range("F4").Value = 1680
$ws1.Range("F5").Value = 2417
$ws1.Range("F8").Value = 1368
$ws1.Range("F9").Value = 1082
$ws1.Range("F14").Value = 94
$ws1.Range("F15").Value = 546
$ws1.Range("F16").Value = 8441
$ws1.Range("F17").Value = 369
$ws1.Range("F18").Value = 2480
$ws1.Range("F19").Value = 247
$ws1.Range("F23").Value = 579
$ws1.Range("F25").Value = 1149
$ws1.Range("F27").Value = 1975
$ws1.Range("F28").Value = 2034
$ws1.Range("F33").Value = 479
$ws1.Range("F35").Value = 37
$ws1.Range("F43").Value = 108
$ws1.Range("F44").Value = 77
$ws1.Range("F45").Value = 252

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 5

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3344
$ws4.Range("G3").Value = 80
$ws4.Range("F4").Value = 1680
$ws4.Range("F5").Value = 2417
$ws4.Range("F8").Value = 1368
$ws4.Range("F10").Value = 1082
$ws4.Range("F14").Value = 94
$ws4.Range("F15").Value = 546
$ws4.Range("F16").Value = 8441
$ws4.Range("F17").Value = 369
$ws4.Range("F18").Value = 2480
$ws4.Range("F20").Value = 247
$ws4.Range("F24").Value = 579
$ws4.Range("F26").Value = 1149
$ws4.Range("F28").Value = 1975
$ws4.Range("F29").Value = 2034
$ws4.Range("F33").Value = 479
$ws4.Range("F35").Value = 37
$ws4.Range("F43").Value = 5
$ws4.Range("F47").Value = 108
$ws4.Range("F49").Value = 252
